# The exported "Saldo" sheet lists one account per row (Conta, Nome, Saldo).
# This change re-inserts the RENATO (account 004862672) record earlier in the
# list with an updated balance, and removes the old entry for that same
# account that used to sit right after the CARLOS (004488571) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row just above the VENIA (004813166) row, which is
# row 38, shifting VENIA and everything below it down by one row.
$ws.Rows.Item(38).Insert()

# Fill the new row with the RENATO record. Force column A to text first so
# the account number keeps its leading zeros instead of being read as a
# number.
$ws.Cells.Item(38, 1).NumberFormat = "@"
$ws.Cells.Item(38, 1).Value = "004862672"
$ws.Cells.Item(38, 2).Value = "RENATO"
$ws.Cells.Item(38, 3).Value = 386.95

# Remove the old RENATO row, which used to immediately follow CARLOS
# (004488571) and, after the insertion above, has shifted down to row 85.
$ws.Rows.Item(85).Delete()
